# The table cell ending in "Tell them whatever they are feeling is okay
# and give them comfort." had a trailing empty paragraph right after it
# (an empty run, no text). Find that text paragraph, then remove the
# very next paragraph, which is empty aside from its paragraph mark
# (and, since it is the last paragraph in the cell, the cell-end mark).

$d = $word.ActiveDocument

$target = "Tell them whatever they are feeling is okay and give them comfort."

$paras = $d.Paragraphs
$count = $paras.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.StartsWith($target)) {
        $nextPara = $paras.Item($i + 1)
        # Strip the paragraph mark (chr 13) and, if present, the
        # end-of-cell mark (chr 7) before checking for emptiness.
        $nextText = $nextPara.Range.Text.TrimEnd([char]13, [char]7)
        if ($nextText -eq "") {
            $nextPara.Range.Delete()
        }
        break
    }
}
